$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '67.779.20'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '3.334.78'
$ws.Range('E3').Value = '  +1.80%  '
$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '581.48'
$cell.Style = $origStyle
$ws.Range('E5').Value = '  +0.88%  '
$cell = $ws.Range('D6')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '176.45'
$cell.Style = $origStyle
$ws.Range('E6').Value = '  +2.43%  '
$ws.Range('E7').Value = '  -0.08%  '
$cell = $ws.Range('D8')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.590'
$cell.Style = $origStyle
$ws.Range('E8').Value = '  +2.18%  '
$ws.Range('D9').Value = '3.330.85'
$ws.Range('E9').Value = '  +1.93%  '
$ws.Range('E10').Value = '  +6.38%  '
$ws.Range('E11').Value = '  +2.13%  '
$cell = $ws.Range('D12')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '47.16'
$cell.Style = $origStyle
$ws.Range('E12').Value = '  +5.01%  '
$ws.Range('E13').Value = '  +0.82%  '
$cell = $ws.Range('D14')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '694.56'
$cell.Style = $origStyle
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('D15').Value = '3.873.88'
$ws.Range('E15').Value = '  +1.97%  '
$cell = $ws.Range('D16')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.44'
$cell.Style = $origStyle
$ws.Range('E16').Value = '  +2.67%  '
$ws.Range('D17').Value = '67.770.84'
$ws.Range('E17').Value = '  +1.31%  '
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').Value = '3.333.27'
$ws.Range('E19').Value = '  +1.80%  '
$cell = $ws.Range('D20')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '17.54'
$cell.Style = $origStyle
$ws.Range('E20').Value = '  +2.19%  '
$cell = $ws.Range('D21')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.07'
$cell.Style = $origStyle
$ws.Range('E21').Value = '  +3.72%  '
$ws.Range('E22').Value = '  +1.27%  '
$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.47'
$cell.Style = $origStyle
$ws.Range('E23').Value = '  +5.30%  '
$cell = $ws.Range('D24')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '17.01'
$cell.Style = $origStyle
$ws.Range('E24').Value = '  +0.73%  '
$cell = $ws.Range('D25')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '100.35'
$cell.Style = $origStyle
$ws.Range('E25').Value = '  +1.34%  '
$cell = $ws.Range('D26')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.91'
$cell.Style = $origStyle
$ws.Range('E26').Value = '  +2.33%  '
$ws.Range('E27').Value = '  +2.47%  '
$cell = $ws.Range('D28')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.55'
$cell.Style = $origStyle
$ws.Range('E28').Value = '  +5.47%  '
$cell = $ws.Range('D29')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '33.06'
$cell.Style = $origStyle
$ws.Range('E29').Value = '  -0.61%  '
$cell = $ws.Range('D30')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.58'
$cell.Style = $origStyle
$ws.Range('E30').Value = '  +3.66%  '
$cell = $ws.Range('D31')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.09'
$cell.Style = $origStyle
$ws.Range('E31').Value = '  +7.59%  '
$cell = $ws.Range('D32')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '566.57'
$cell.Style = $origStyle
$ws.Range('E32').Value = '  -1.93%  '
$cell = $ws.Range('D33')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.01'
$cell.Style = $origStyle
$ws.Range('E33').Value = '  +1.81%  '
$cell = $ws.Range('D34')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.105'
$cell.Style = $origStyle
$ws.Range('E34').Value = '  +3.61%  '
$cell = $ws.Range('D35')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '57.34'
$cell.Style = $origStyle
$ws.Range('E35').Value = '  +4.21%  '
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').Value = '3.709.67'
$ws.Range('E37').Value = '  -2.66%  '
$ws.Range('E38').Value = '  +2.35%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Range('D39')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.134'
$cell.Style = $origStyle
$ws.Range('E39').Value = '  +4.69%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range('D40')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '34.72'
$cell.Style = $origStyle
$ws.Range('E40').Value = '  +10.71%  '
$ws.Range('E41').Value = '  +3.29%  '
$cell = $ws.Range('D42')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.16'
$cell.Style = $origStyle
$ws.Range('E42').Value = '  +6.65%  '
$ws.Range('D43').Value = '0.0₃0675'
$ws.Range('E43').Value = '  +2.33%  '
$ws.Range('E44').Value = '  +4.14%  '
$cell = $ws.Range('D45')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.29'
$cell.Style = $origStyle
$ws.Range('E45').Value = '  -2.40%  '
$ws.Range('E46').Value = '  +2.17%  '
$ws.Range('E47').Value = '  +5.42%  '
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('E49').Value = '  +0.04%  '
$cell = $ws.Range('D50')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.34'
$cell.Style = $origStyle
$ws.Range('E50').Value = '  -0.24%  '
$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '131.21'
$cell.Style = $origStyle
$ws.Range('E51').Value = '  +1.79%  '
